$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency price/volume figures (refresh snapshot).
# Some "Price" (column D) values parse as plain numbers (e.g. "7.390"), but the
# source data is textual (trailing zeros / grouped digits must be preserved), so
# those cells are explicitly formatted as Text before assignment to stop Excel
# from silently converting them to numeric values and dropping significant digits.

$ws.Range("D2").Value = "27.493.36"
$ws.Range("E2").Value = "  +5.43%  "
$ws.Range("D3").Value = "1.725.27"
$ws.Range("E3").Value = "  +4.72%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.38"
$ws.Range("E5").Value = "  +3.24%  "
$ws.Range("E6").Value = "  +3.11%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2665"
$ws.Range("E8").Value = "  +1.56%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06583"
$ws.Range("E9").Value = "  +4.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.57"
$ws.Range("E10").Value = "  +6.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07696"
$ws.Range("E11").Value = "  +0.17%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.601"
$ws.Range("E12").Value = "  +0.27%  "
$ws.Range("D13").Value = "1.727.60"
$ws.Range("E13").Value = "  +7.53%  "
$ws.Range("D14").Value = "1.963.13"
$ws.Range("E14").Value = "  +4.69%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5814"
$ws.Range("E15").Value = "  +4.49%  "
$ws.Range("D16").Value = "0.0₅8269"
$ws.Range("E16").Value = "  +2.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.77"
$ws.Range("E17").Value = "  +4.17%  "
$ws.Range("D18").Value = "27.503.40"
$ws.Range("E18").Value = "  +5.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.83"
$ws.Range("E19").Value = "  +13.05%  "
$ws.Range("E20").Value = "  +0.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.719"
$ws.Range("E21").Value = "  +2.49%  "
$ws.Range("E22").Value = "  +1.63%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.073"
$ws.Range("E23").Value = "  +2.74%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.30"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.752"
$ws.Range("E26").Value = "  +16.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1233"
$ws.Range("E27").Value = "  +4.44%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.390"
$ws.Range("E28").Value = "  +2.93%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.51"
$ws.Range("E29").Value = "  +4.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05484"
$ws.Range("E30").Value = "  +2.76%  "
$ws.Range("E31").Value = "  +2.67%  "
$ws.Range("E33").Value = "  +3.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.657"
$ws.Range("E34").Value = "  +7.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.862"
$ws.Range("E35").Value = "  +2.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9632"
$ws.Range("E36").Value = "  +2.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.424"
$ws.Range("E37").Value = "  +0.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5966"
$ws.Range("E38").Value = "  +6.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01648"
$ws.Range("E39").Value = "  +4.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.904"
$ws.Range("E40").Value = "  +2.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8522"
$ws.Range("E41").Value = "  +3.44%  "
$ws.Range("D42").Value = "1.053.22"
$ws.Range("E42").Value = "  +2.72%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("E44").Value = "  +0.46%  "
$ws.Range("D45").Value = "1.869.53"
$ws.Range("E45").Value = "  +4.68%  "
$ws.Range("E46").Value = "  +3.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "58.81"
$ws.Range("E47").Value = "  +2.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4474"
$ws.Range("E48").Value = "  +3.82%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.170"
$ws.Range("E49").Value = "  +3.49%  "
$ws.Range("E50").Value = "  -0.10%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05237"
$ws.Range("E51").Value = "  +2.52%  "
